$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(8, 1, 10536,    2025, "Bibi Cell Mundi"),
    @(8, 2, 10676.75, 2025, "Bibi Cell Manauara"),
    @(8, 3, 16076,    2025, "Bibi Cell Vieiralves"),
    @(8, 4, 10218.55, 2025, "Bibi Cell Ponta Negra")
)

$startRow = 264
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}
